# edit.ps1 - reproduces the OOXML diff:
#   * removes the old "{{PASSWORD}}" textbox (id 23 "TextBox 23")
#   * turns the old "{{QR_WIFI}}" textbox (id 29 "TextBox 28") into a brand
#     new textbox holding the password text (renamed "TextBox 23"), moved to
#     the top of the block, with the Circe font applied to the run
#   * appends a further brand new textbox holding the "{{QR_WIFI}}" text
#     (named "TextBox 30"), with wrap="square" instead of wrap="none"
#
# PowerPoint stores shape Left/Top/Width/Height as single-precision (float32)
# point values, and persists them back to EMU (1 pt = 12700 EMU) by
# truncating - exactly like real PowerPoint COM automation. To land on an
# exact target EMU value we search for the nearest point value whose
# float32 round-trip truncates to that EMU.
function ConvertTo-Emu($pt) {
    $f = [single]$pt
    return [math]::Floor([double]$f * 12700)
}

function ConvertFrom-Emu($targetEmu) {
    $pt = $targetEmu / 12700.0
    $step = 0.0000001
    for ($i = 0; $i -lt 200000; $i++) {
        $e = ConvertTo-Emu $pt
        if ($e -eq $targetEmu) {
            return $pt
        } elseif ($e -lt $targetEmu) {
            $pt = $pt + $step
        } else {
            $pt = $pt - $step
        }
    }
    return $pt
}

function Find-ShapeByText($slide, $matchText) {
    $n = $slide.Shapes.Count
    for ($i = 1; $i -le $n; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $matchText) {
                return $sh
            }
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# Locate the slide that contains the password / QR placeholders.
$slide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $candidate = $p.Slides.Item($si)
    if ((Find-ShapeByText $candidate "{{PASSWORD}}") -ne $null) {
        $slide = $candidate
        break
    }
}

# Remember the old QR textbox's geometry before we delete it.
$qrShape = Find-ShapeByText $slide "{{QR_WIFI}}"
$qrLeftPt = $qrShape.Left
$qrTopPt = $qrShape.Top
$qrWidthPt = $qrShape.Width
$qrHeightPt = $qrShape.Height

# 1) Delete the old password textbox (id 23 "TextBox 23").
$pwShape = Find-ShapeByText $slide "{{PASSWORD}}"
$pwShape.Delete()

# 2) Delete the old QR textbox (id 29, also named "TextBox 28").
$qrShape.Delete()

# 3) Add the new password textbox (becomes id 30 "TextBox 23").
$newPwLeft = ConvertFrom-Emu 888976
$newPwTop = ConvertFrom-Emu 3565121
$newPwWidth = ConvertFrom-Emu 2350430
$newPwHeight = ConvertFrom-Emu 396262

$newPw = $slide.Shapes.AddTextbox(1, $newPwLeft, $newPwTop, $newPwWidth, $newPwHeight)
$newPw.Name = "TextBox 23"

$pwFrame = $newPw.TextFrame
$pwFrame.MarginLeft = 0
$pwFrame.MarginTop = 0
$pwFrame.MarginRight = 0
$pwFrame.MarginBottom = 0
$pwFrame.VerticalAnchor = 1
$pwFrame.AutoSize = 1

$pwRange = $pwFrame.TextRange
$pwRange.Text = "{{PASSWORD}}"
$pwRange.Font.Name = "Circe"
$pwRange.ParagraphFormat.Alignment = 2
$pwRange.ParagraphFormat.SpaceWithin = 33.5

# Re-apply the exact geometry last, since AutoSize/spAutoFit recalculates
# the height as soon as the text/font is set.
$newPw.Left = $newPwLeft
$newPw.Top = $newPwTop
$newPw.Width = $newPwWidth
$newPw.Height = $newPwHeight

# 4) Add the new QR textbox (becomes id 31 "TextBox 30"), reusing the old
#    QR shape's geometry family (shifted per the diff) and text.
$newQrLeft = ConvertFrom-Emu 1143198
$newQrTop = ConvertFrom-Emu 4904003
$newQrWidth = ConvertFrom-Emu 1307902
$newQrHeight = ConvertFrom-Emu 369332

$newQr = $slide.Shapes.AddTextbox(1, $newQrLeft, $newQrTop, $newQrWidth, $newQrHeight)
# Default auto-generated name is already "TextBox 30"; set explicitly to be safe.
$newQr.Name = "TextBox 30"
$newQr.Fill.Visible = 0

$qrFrame = $newQr.TextFrame
$qrFrame.WordWrap = -1
$qrFrame.AutoSize = 1

$qrRange = $qrFrame.TextRange
$qrRange.Text = "{{QR_WIFI}}"

$newQr.Left = $newQrLeft
$newQr.Top = $newQrTop
$newQr.Width = $newQrWidth
$newQr.Height = $newQrHeight
